$wb = $excel.ActiveWorkbook

$postSheet = $wb.Worksheets.Item("PostData")
$postSheet.Activate()
$postSheet.Range("G6").Select()

$paySheet = $wb.Worksheets.Item("Payment")
$paySheet.Activate()
$paySheet.Range("L2").Value = "Online"
$paySheet.Range("K7").Select()
